# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new first sheet "Player Info" with ID/NAME/BATTING_HAND/BOWL_STYLE.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on both "ODI Batting" and "ODI Bowling",
#    replacing the full scorecard URL values with the bare numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update MATCH_CARD_LINK -> MATCH_CODE on the two existing sheets
# ---------------------------------------------------------------------
function Update-MatchCodeColumn($ws, $col) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Value2 = "MATCH_CODE"

    $row = 2
    while ($true) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value2
        if ($val -eq $null -or $val -eq "") { break }
        if ($val -match "MatchCode=(\d+)") {
            # keep the cell as text so the bare numeric match code (e.g.
            # "3987") is stored the same way the URL string used to be,
            # not converted into a number
            $cell.NumberFormat = "@"
            $cell.Value2 = $matches[1]
        }
        $row = $row + 1
    }
}

$battingSheet = $wb.Worksheets.Item("ODI Batting")
Update-MatchCodeColumn $battingSheet 4   # column D

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
Update-MatchCodeColumn $bowlingSheet 2   # column B

# ---------------------------------------------------------------------
# Step 2: insert the new "Player Info" sheet as the first sheet
# ---------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add()
$infoSheet.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $infoSheet.Cells.Item(1, $i + 1)
    $cell.Value2 = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$values = @("4644", "Lahiru Dilshan Madushanka", "Right Handed", "Right Arm Fast Medium")
for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $infoSheet.Cells.Item(2, $i + 1)
    if ($i -eq 0) {
        # column A holds the numeric-looking player ID ("4644") as text,
        # like the rest of this workbook stores its ID/code columns
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $values[$i]
}

$infoSheet.Range("A1").Select()
